$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 56 ("cwl_warn_processor") - the id cell A56 was still using the
#    "log" (default/black) font color instead of the "warn" (orange) one
#    used by every other cwl_warn_* row. Fix the font color to match the
#    existing orange accent (RGB FFC000 -> OLE BGR 49407).
$ws.Cells.Item(56, 1).Font.Color = 49407

# 2. New row 57 entry: "cwl_log_ele_gain" / "auto gained ability id: {0} on {1}"
#    First copy the formatting from row 56's C/D cells (wrapped, same font)
#    onto row 57's C/D cells so the new row matches the sheet's look.
$ws.Range("C56:D56").Copy()
$ws.Range("C57:D57").PasteSpecial(-4122)

$ws.Cells.Item(57, 1).Value = "cwl_log_ele_gain"
$ws.Cells.Item(57, 3).Value = "auto gained ability id: {0} on {1}"
$ws.Cells.Item(57, 4).Value = "auto gained ability id: {0} on {1}"

# 3. The active selection moved from D58 to D59.
$ws.Range("D59").Select()
